# DCL_ENTC_TYB.xlsx — add the 2025-06-23 attendance column.
#
# The attendance sheet gained one more date column (U) on the right of
# "Attendance %": a new header "2025-06-23", a "❌" mark for every student
# row, and the per-row "Total" (col S) bumped by one to account for the
# extra day being counted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- U1: header cell, same look as the other date headers (bold / bordered
# / centered, i.e. the style already used by T1) but literal text so Excel
# doesn't reinterpret "2025-06-23" as a serial date.
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)          # xlPasteFormats: clone T1's style only
$ws.Range("U1").Formula = "=""2025-06-23"""  # literal-text formula avoids date coercion
$ws.Range("U1").Copy()
$ws.Range("U1").PasteSpecial(-4163)          # xlPasteValues: collapse formula -> plain value, keep style

# --- U2 / U3: attendance mark for each student on the new date.
$ws.Range("U2").Value = "❌"
$ws.Range("U3").Value = "❌"

# --- S2 / S3: "Total" count grows by one with the extra date column.
$ws.Range("S2").Value = 16
$ws.Range("S3").Value = 16

$excel.CutCopyMode = $false
